$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.299.83'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.594.41'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '189.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.592.28'
$ws.Range('E7').Value = '  -0.92%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.630'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.661'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000311'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '4.182.21'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '19.77'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '3.602.53'
$ws.Range('D18').Value = '70.318.14'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '489.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '19.59'
$ws.Range('E23').Value = '  +1.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -9.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.07%  '
$ws.Range('D28').Value = '10.99'
$ws.Range('E28').Value = '  -3.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.52%  '
$ws.Range('D32').Value = '12.25'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '66.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  -2.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '575.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('D37').Value = '0.0₃0808'
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.398'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.16%  '
$ws.Range('D41').Value = '3.23'
$ws.Range('E41').Value = '  +16.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('E43').Value = '  -6.42%  '
$ws.Range('D44').Value = '3.212.92'
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').Value = '3.04'
$ws.Range('E45').Value = '  -3.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0446'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.04%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.94%  '
